$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 52-61 (zz033 "lydata" entries added)
# Columns: A=lya B=snöfri procent C=snöfri area (m^2) D=lufttemp
#          E=marktemperatur orange F=marktemperatur svart G=underlag marktemp
#          H=snödjup I=riktning (grader) J=vinkel K=antal lyöppningar L=aktiv
$rows = @(
    @{ r=52; B=75; C="91.5"; D=0; E=-0.5; F=-0.6; G="s";  H=17; I="250 SV"; J=20; K=3; L="j" },
    @{ r=53; B=75; C="91.5"; D=0; E=-0.2; F=-0.3; G="ns"; H=12; I="250 SV"; J=20; K=3; L="j" },
    @{ r=54; B=75; C="91.5"; D=0; E=3.8;  F=4;    G="b";  H=0;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=55; B=75; C="91.5"; D=0; E=-0.1; F=-0.2; G="ns"; H=8;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=56; B=75; C="91.5"; D=0; E=-0.2; F=-0.1; G="b";  H=0;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=57; B=75; C="91.5"; D=0; E=0.3;  F=0.1;  G="b";  H=0;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=58; B=75; C="91.5"; D=0; E=-0.4; F=-0.5; G="ns"; H=22; I="250 SV"; J=20; K=3; L="j" },
    @{ r=59; B=75; C="91.5"; D=0; E=-0.3; F=-0.5; G="ns"; H=5;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=60; B=75; C="91.5"; D=0; E=-0.2; F=-0.3; G="ns"; H=3;  I="250 SV"; J=20; K=3; L="j" },
    @{ r=61; B=75; C="91.5"; D=0; E=-0.4; F=-0.4; G="b";  H=0;  I="250 SV"; J=20; K=3; L="j" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 2).Value = $row.B
    # Column C holds a numeric-looking value ("91.5") that must be stored as text
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
}

# Update frozen pane top-left cell and active selection to reflect scrolled view
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("G62").Select()
